$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-11 hold a cyclic rotation of the Id (A), Antal (I), Ost (Q) and
# Nord (R) columns: row9 <- old row10, row10 <- old row11, row11 <- old row9.

# Row 9
$ws.Range("A9").Value = 111675585
$ws.Range("I9").Value = "'1"
$ws.Range("I9").Style = "Normal"
$ws.Range("Q9").Value = 690349.9096738817
$ws.Range("R9").Value = 6661440.004307052

# Row 10
$ws.Range("A10").Value = 111675586
$ws.Range("I10").Value = "'2"
$ws.Range("I10").Style = "Normal"
$ws.Range("Q10").Value = 690348.8581766916
$ws.Range("R10").Value = 6661440.95072202

# Row 11
$ws.Range("A11").Value = 111675587
$ws.Range("I11").Value = "'3"
$ws.Range("I11").Style = "Normal"
$ws.Range("Q11").Value = 690344.8588249951
$ws.Range("R11").Value = 6661440.743740954
